# Applies the periodic "cryptos list" price/volume refresh described by the commit.
# Only the cells that actually changed are touched; everything else is left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '28.344.22'
$ws.Range('E2').Value = '  -0.85%  '

# Row 3
$ws.Range('D3').Value = '1.551.08'
$ws.Range('E3').Value = '  -1.63%  '

# Row 5
$ws.Range('D5').Value = '''210.09'
$ws.Range('E5').Value = '  -1.51%  '

# Row 6
$ws.Range('D6').Value = '''0.480'
$ws.Range('E6').Value = '  -2.03%  '

# Row 7
$ws.Range('E7').Value = '  -0.25%  '

# Row 8
$ws.Range('D8').Value = '''23.88'
$ws.Range('E8').Value = '  -0.54%  '

# Row 9
$ws.Range('E9').Value = '  -1.90%  '

# Row 10
$ws.Range('D10').Value = '''0.0582'
$ws.Range('E10').Value = '  -1.47%  '

# Row 11
$ws.Range('D11').Value = '''0.0888'
$ws.Range('E11').Value = '  -0.26%  '

# Row 12
$ws.Range('D12').Value = '1.774.44'
$ws.Range('E12').Value = '  -1.50%  '

# Row 13
$ws.Range('D13').Value = '1.552.28'
$ws.Range('E13').Value = '  -1.56%  '

# Row 14
$ws.Range('D14').Value = '28.348.60'
$ws.Range('E14').Value = '  -0.78%  '

# Row 15
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '''3.62'
$ws.Range('E15').Value = '  -1.76%  '

# Row 16
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '''0.510'
$ws.Range('E16').Value = '  -1.83%  '

# Row 17
$ws.Range('D17').Value = '''60.82'
$ws.Range('E17').Value = '  -2.09%  '

# Row 18
$ws.Range('D18').Value = '''227.62'
$ws.Range('E18').Value = '  -1.24%  '

# Row 19
$ws.Range('E19').Value = '  -0.44%  '

# Row 20
$ws.Range('E20').Value = '  -2.28%  '

# Row 21
$ws.Range('E21').Value = '  -0.22%  '

# Row 22
$ws.Range('E22').Value = '  +0.72%  '

# Row 23
$ws.Range('D23').Value = '''8.93'
$ws.Range('E23').Value = '  -2.55%  '

# Row 24
$ws.Range('D24').Value = '''2.03'
$ws.Range('E24').Value = '  -1.60%  '

# Row 25
$ws.Range('D25').Value = '''150.86'
$ws.Range('E25').Value = '  -0.27%  '

# Row 26
$ws.Range('D26').Value = '''14.72'
$ws.Range('E26').Value = '  -1.81%  '

# Row 27
$ws.Range('E27').Value = '  -1.19%  '

# Row 29
$ws.Range('D29').Value = '''6.24'
$ws.Range('E29').Value = '  -3.10%  '

# Row 30
$ws.Range('D30').Value = '''0.0467'
$ws.Range('E30').Value = '  -3.10%  '

# Row 31
$ws.Range('E31').Value = '  -4.37%  '

# Row 32
$ws.Range('D32').Value = '''3.16'
$ws.Range('E32').Value = '  -1.27%  '

# Row 33
$ws.Range('D33').Value = '1.386.12'
$ws.Range('E33').Value = '  -0.80%  '

# Row 34
$ws.Range('D34').Value = '''3.01'
$ws.Range('E34').Value = '  -2.92%  '

# Row 35
$ws.Range('D35').Value = '''1.08'
$ws.Range('E35').Value = '  +2.04%  '

# Row 36
$ws.Range('E36').Value = '  -3.54%  '

# Row 37
$ws.Range('E37').Value = '  -1.09%  '

# Row 38
$ws.Range('E38').Value = '  -1.69%  '

# Row 39
$ws.Range('E39').Value = '  -2.98%  '

# Row 40
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''1.92'
$ws.Range('E40').Value = '  +2.04%  '

# Row 41
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').Value = '''0.511'
$ws.Range('E41').Value = '  -1.90%  '

# Row 42
$ws.Range('E42').Value = '  -0.28%  '

# Row 43
$ws.Range('D43').Value = '''0.777'
$ws.Range('E43').Value = '  -2.06%  '

# Row 44
$ws.Range('D44').Value = '''0.0454'
$ws.Range('E44').Value = '  -2.10%  '

# Row 46
$ws.Range('E46').Value = '  -1.62%  '

# Row 47
$ws.Range('D47').Value = '1.686.91'
$ws.Range('E47').Value = '  -1.58%  '

# Row 48
$ws.Range('D48').Value = '''0.861'
$ws.Range('E48').Value = '  -10.47%  '

# Row 49
$ws.Range('D49').Value = '''85.32'
$ws.Range('E49').Value = '  -1.39%  '

# Row 50
$ws.Range('D50').Value = '''42.53'
$ws.Range('E50').Value = '  +6.56%  '

# Row 51
$ws.Range('E51').Value = '  +0.28%  '
